$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (cellRef, newValue, isNumericLike) triples describing every cell touched by the commit.
# Numeric-looking strings are written via a literal-text formula then pasted back as
# values-only, so Excel keeps them as text (matching the source inlineStr cells) instead
# of auto-converting to a Number/Percentage cell.
$edits = @(
    @("D2", "311.57"),
    @("E2", "1.30%"),
    @("D3", "37.76"),
    @("E3", "0.38%"),
    @("D4", "5.118"),
    @("E4", "0.32%"),
    @("D5", "0.07909"),
    @("E5", "0.34%"),
    @("B6", "GateToken"),
    @("C6", "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"),
    @("D6", "4.403"),
    @("E6", "1.58%"),
    @("B7", "FTXToken"),
    @("C7", "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"),
    @("D7", "1.907"),
    @("E7", "-3.69%"),
    @("B8", "KuCoinToken"),
    @("C8", "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"),
    @("D8", "8.237"),
    @("E8", "0.03%"),
    @("D9", "0.9271"),
    @("E9", "-0.40%"),
    @("D10", "0.1201"),
    @("E10", "-7.72%"),
    @("D11", "0.1911"),
    @("E11", "0.66%"),
    @("D12", "0.09323"),
    @("E12", "4.98%"),
    @("D13", "0.03382"),
    @("E13", "-1.59%"),
    @("D14", "0.09614"),
    @("E14", "-1.33%"),
    @("D15", "0.001373"),
    @("E15", "-1.30%"),
    @("D16", "0.005868"),
    @("E16", "-0.15%"),
    @("D17", "3.530"),
    @("E17", "-0.93%"),
    @("B18", "BTSEToken"),
    @("C18", "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"),
    @("D18", "3.079"),
    @("E18", "-0.70%"),
    @("D19", "0.3449"),
    @("E19", "0.53%"),
    @("D20", "5.269"),
    @("E20", "5.37%"),
    @("E21", "-0.50%"),
    @("E23", "179.62%"),
    @("E24", "1.27%"),
    @("D25", "0.001249"),
    @("E25", "2.46%"),
    @("D26", "0.004274"),
    @("E26", "-7.26%"),
    @("D27", "0.0001297"),
    @("E27", "-63.91%"),
    @("D39", "0.02087"),
    @("E39", "-10.89%"),
    @("D40", "0.05069"),
    @("E40", "0.70%"),
    @("D41", "0.007612"),
    @("E41", "1.05%"),
    @("D42", "0.009106"),
    @("E42", "-6.62%"),
    @("D43", "0.1353"),
    @("E43", "-0.13%"),
    @("D44", "0.002085"),
    @("E44", "-0.24%"),
    @("D45", "0.008629"),
    @("E45", "7.73%"),
    @("D46", "0.00006683"),
    @("E46", "2.45%"),
    @("D47", "0.00000000748"),
    @("E47", "-0.21%"),
    @("D48", "0.002882"),
    @("E48", "-3.84%"),
    @("D49", "0.001197"),
    @("E49", "-0.24%"),
    @("D50", "0.00002095"),
    @("E50", "-0.21%"),
    @("D51", "0.0001995"),
    @("E51", "-0.21%"),
)

foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newVal = $edit[1]
    $range = $ws.Range($cellRef)
    # Build a formula whose result is the literal text (escape any embedded quotes for the formula string)
    $escaped = $newVal.Replace("""", """""")
    $range.Formula = "=""" + $escaped + """"
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0
